$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.47 = 39090.91 pesos`n✅ 39090.91 pesos = 9.4 = 958.48 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 4128
$ws2.Range("N12").Value = 4160
$ws2.Range("O12").Value = 102
